$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: __init__ row gains a Range column entry ---
$ws.Range("G2").Value = "(minimum, maximum, dice_count)"

# --- Row 3: __str__ ---
$ws.Range("G3").Value = "(verbose), str"

# --- Row 4: __unicode__ ---
$ws.Range("G4").Value = "inherit(Die)"

# --- Row 5: valid ---
$ws.Range("G5").Value = "@property(), bool"

# --- Row 6: rolled ---
$ws.Range("G6").Value = "@property(), bool"

# --- Row 7: sides ---
$ws.Range("G7").Value = "@property(), int"

# --- Row 8: average ---
$ws.Range("G8").Value = "@property(), num"

# --- Row 9: result ---
$ws.Range("G9").Value = "@property(), int"

# --- Row 10: minimum ---
$ws.Range("G10").Value = "inherit(Die)"

# --- Row 11: maximum -- Range no longer has its own maximum setter, now inherits Die;
#     D/E/F (Pecentile/D1000/D10000) also now simply inherit(Die) instead of bespoke text ---
$ws.Range("D11").Value = "inherit(Die)"
$ws.Range("E11").Value = "inherit(Die)"
$ws.Range("F11").Value = "inherit(Die)"
$ws.Range("G11").Value = "@property(), int"

# --- Row 12: roll ---
$ws.Range("G12").Value = "(), None"

# --- Row 13: history ---
$ws.Range("G13").Value = "inherit(Die)"

# --- Row 14: clear_history ---
$ws.Range("G14").Value = "inherit(Die)"

# --- Row 15: dice ---
$ws.Range("G15").Value = "@property(), list[Die]"

# --- Row 16: ones ---
$ws.Range("G16").Value = "Not Implemented"

# --- Row 17: tens ---
$ws.Range("G17").Value = "Not Implemented"

# --- Row 18: hundreds ---
$ws.Range("G18").Value = "Not Implemented"

# --- Row 19: thousands ---
$ws.Range("G19").Value = "Not Implemented"

# --- New row 20: dice_count (removed as a Range setter; now just a method entry) ---
$ws.Range("A20").Value = "dice_count"
$ws.Range("B20").Value = "Not Implemented"
$ws.Range("C20").Value = "Not Implemented"
$ws.Range("D20").Value = "Not Implemented"
$ws.Range("E20").Value = "Not Implemented"
$ws.Range("F20").Value = "Not Implemented"
$ws.Range("G20").Value = "@property(), int"

# --- New row 21: dice_sort ---
$ws.Range("A21").Value = "dice_sort"
$ws.Range("B21").Value = "Not Implemented"
$ws.Range("C21").Value = "Not Implemented"
$ws.Range("D21").Value = "Not Implemented"
$ws.Range("E21").Value = "Not Implemented"
$ws.Range("F21").Value = "Not Implemented"
$ws.Range("G21").Value = "(), None"

# --- New row 22: _build_total ---
$ws.Range("A22").Value = "_build_total"
$ws.Range("B22").Value = "Not Implemented"
$ws.Range("C22").Value = "Not Implemented"
$ws.Range("D22").Value = "Not Implemented"
$ws.Range("E22").Value = "Not Implemented"
$ws.Range("F22").Value = "Not Implemented"
$ws.Range("G22").Value = "(dice_to_exclude), None"

# --- Column G widened to fit the longer Range descriptions ---
$ws.Columns.Item(7).ColumnWidth = 31.6

# --- Selection moved ---
$ws.Range("C29").Select()
